$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-20 Wednesday" "2023-09-21 Thursday"

Replace-Text "60÷9=6, 6" "41÷9=4, 5"
Replace-Text "58÷3=19, 1" "16÷3=5, 1"
Replace-Text "28÷4=7, 0" "45÷4=11, 1"
Replace-Text "54÷2=27, 0" "69÷8=8, 5"
Replace-Text "46÷2=23, 0" "69÷4=17, 1"

Replace-Text "86÷5=17, 1" "17÷9=1, 8"
Replace-Text "59÷9=6, 5" "86÷8=10, 6"
Replace-Text "94÷6=15, 4" "19÷2=9, 1"
Replace-Text "12÷5=2, 2" "82÷3=27, 1"
Replace-Text "41÷8=5, 1" "37÷4=9, 1"

Replace-Text "33÷4=8, 1" "72÷5=14, 2"
Replace-Text "56÷6=9, 2" "52÷7=7, 3"
Replace-Text "15÷8=1, 7" "97÷9=10, 7"
Replace-Text "52÷8=6, 4" "13÷2=6, 1"
Replace-Text "23÷4=5, 3" "50÷9=5, 5"

Replace-Text "84÷7=12, 0" "21÷2=10, 1"
Replace-Text "18÷2=9, 0" "87÷4=21, 3"
Replace-Text "62÷6=10, 2" "56÷5=11, 1"
Replace-Text "32÷4=8, 0" "12÷2=6, 0"
Replace-Text "86÷3=28, 2" "47÷7=6, 5"

Replace-Text "44÷5=8, 4" "75÷3=25, 0"
Replace-Text "44÷2=22, 0" "34÷4=8, 2"
Replace-Text "66÷3=22, 0" "95÷7=13, 4"
Replace-Text "30÷8=3, 6" "55÷7=7, 6"
Replace-Text "42÷2=21, 0" "30÷4=7, 2"
